# Generate Report for Handback
# Refresh the handback-status report with the latest generation/handoff/
# handback timestamps and correct the zh-cn / de-de "Priority" value
# (ht -> mt) for the two in-progress rows.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" for the two rows that
#     just finished a fresh handback pass.
$wsOverview.Range("G2").Value = "2016-08-26 00:16:31"
$wsOverview.Range("G3").Value = "2016-08-26 00:16:31"

# --- zh-cn sheet
# Priority flips from "ht" (human translation) to "mt" (machine translation)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-26 00:16:27"
$wsZhCn.Range("H3").Value = "2016-08-26 00:16:27"
# Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-26 00:16:42"
$wsZhCn.Range("K3").Value = "2016-08-26 00:16:42"

# --- de-de sheet
# Priority flips from "ht" (human translation) to "mt" (machine translation)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# Correspond Handoff Datetime (mirrors Overview's Latest HO Xliff Generate Date)
$wsDeDe.Range("H2").Value = "2016-08-26 00:16:31"
$wsDeDe.Range("H3").Value = "2016-08-26 00:16:31"
# Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-26 00:16:49"
$wsDeDe.Range("K3").Value = "2016-08-26 00:16:49"
